# Slide 17 ("Method checkConstraints() for Class Variable") - Content
# Placeholder 2, second paragraph: split the single run
#   "type T is array[10] of Integer;"
# into three runs so that "is" becomes "=":
#   "type " + "T = " + "array[10] of Integer;"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(17)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# The paragraph that currently reads "type T is array[10] of Integer;"
$declPara = $tr.Paragraphs(2, 1)

# Replace the "T is " substring (characters 6-10) with "T = ",
# which splits the original single run into the three runs seen
# in the target: "type " / "T = " / "array[10] of Integer;"
$mid = $declPara.Characters(6, 5)
$mid.Text = "T = "
